$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove discontinued components (bottom-up so row numbers stay valid)
#    row 9 = Netzteil, row 6 = Mikrofon, row 4 = Infrarotsensor
# ------------------------------------------------------------------
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()

# Rows are now:
#  1 header
#  2 Mikrokontroller | 3 Kamera mit Infrarot | 4 Display | 5 LED-Leuchten
#  6 Speicher -> Spiegelfolie | 7 Verbindungskabel... -> Glasscheibe
#  8 Gehaeuse/Rahmen (unchanged) | 9 Spiegel -> Verbindungskabel und Stecker (total-style row)
#  10 Gesamt

# ------------------------------------------------------------------
# 2. Rename remaining component labels to match the new plan
# ------------------------------------------------------------------
$ws.Range("B6").Value = "Spiegelfolie"
$ws.Range("B7").Value = "Glasscheibe"
$ws.Range("B9").Value = "Verbindungskabel und Stecker"

# ------------------------------------------------------------------
# 3. Header row: new columns Versand / Gesamtkosten / Link
# ------------------------------------------------------------------
$ws.Range("D1").Value = "Versand:"
$ws.Range("E1").Value = "Gesamtkosten"
$ws.Range("F1").Value = "Link:"

# ------------------------------------------------------------------
# 4. Apply the existing Euro number style (same as column C) to the new
#    D/E columns before filling in numbers/formulas, so blank cells keep
#    the right formatting too.
# ------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("D2:D8").PasteSpecial(-4122)
$ws.Range("E2:E8").PasteSpecial(-4122)
$ws.Range("D10:E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 9 (bottom-bordered total-style row) uses the bordered Euro style,
# matching the existing C9 cell.
$ws.Range("C9").Copy()
$ws.Range("D9:E9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# F9 only needs the plain bottom-border style (same as A9/B9).
$ws.Range("A9").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 5. Costs (Kosten) + shipping (Versand) per component
# ------------------------------------------------------------------
$ws.Range("C2").Value = 46.99
$ws.Range("D2").Value = 5.99

$ws.Range("C3").Value = 28.59
$ws.Range("D3").Value = 5.95

$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()

$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()

$ws.Range("C6").Value = 14.99
$ws.Range("D6").Value = 0

$ws.Range("C7").Value = 24.99
$ws.Range("D7").Value = 0

$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()

$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()

# ------------------------------------------------------------------
# 6. Gesamtkosten (Kosten + Versand) per row, and grand totals
# ------------------------------------------------------------------
$ws.Range("E2:E9").Formula = "=SUM(C2,D2)"

$ws.Range("C10").Formula = "=SUM(C2:C9)"
$ws.Range("D10").Formula = "=SUM(D2:D9)"
$ws.Range("E10").Formula = "=SUM(E2:E9)"

# ------------------------------------------------------------------
# 7. Links to the shops each component was bought from
# ------------------------------------------------------------------
$ws.Range("F2").Value = "rasppishop"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.rasppishop.de/")

$ws.Range("F3").Value = "buyzero"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://buyzero.de/")

$ws.Range("F6").Value = "Amazon"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.amazon.de/")

$ws.Range("F7").Value = "Amazon"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.amazon.de/")

# F4 carries the hyperlink-cell look even though it has no link of its own.
$ws.Range("F2").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 8. Column widths (best-fit for the new content/columns)
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 28.5703125
$ws.Columns.Item(3).ColumnWidth = 8.140625
$ws.Columns.Item(4).ColumnWidth = 8.85546875
$ws.Columns.Item(5).ColumnWidth = 13.85546875
$ws.Columns.Item(6).ColumnWidth = 10.7109375

# ------------------------------------------------------------------
# 9. Row heights for the total-style row & the grand-total row
# ------------------------------------------------------------------
$ws.Rows.Item(9).RowHeight = 15.75
$ws.Rows.Item(10).RowHeight = 15.75

# ------------------------------------------------------------------
# 10. View tweaks
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 190
$ws.Range("D14").Select()

$wb.Save()
